# New5VIN_CA_SELECT.xlsx - "Add unique VINs to Each CA Select Test PT2" edit
#
# 1. Change the VIN value (shared string) used by A2:A5 from "EEENK3CC&F" to
#    "EEENK2CC&F". Setting the same new value into every cell that shared the
#    original string lets the workbook de-duplicate them back into a single
#    shared-string entry, exactly as the source diff shows (the <si> text is
#    edited in place rather than a new entry being appended).
# 2. Move the sheet's active selection from B10 to B14.
# 3. Narrow column A from width 25 to ~20.71 characters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the VIN column (A2:A5) to the new shared string value.
$ws.Range("A2").Value = "EEENK2CC&F"
$ws.Range("A3").Value = "EEENK2CC&F"
$ws.Range("A4").Value = "EEENK2CC&F"
$ws.Range("A5").Value = "EEENK2CC&F"

# 2) Update the selected/active cell shown in the sheet view.
$ws.Range("B14").Select()

# 3) Resize column A.
$ws.Columns.Item(1).ColumnWidth = 19.833333333333332
